# feat: melhorias e alterações
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: B3, C3, E3, F3, G3 become numeric values instead of text
$ws.Cells.Item(3, 2).Value = 2
$ws.Cells.Item(3, 3).Value = 3
$ws.Cells.Item(3, 5).Value = 5
$ws.Cells.Item(3, 6).Value = 6
$ws.Cells.Item(3, 7).Value = 4

# Row 4: mostly blank row, only D4 (date) and H4 (status) populated
$ws.Cells.Item(4, 4).Value = "2025-01-30 11:22:54"
$ws.Cells.Item(4, 8).Value = "em dia"

# Row 5: fully populated new row (numeric-looking values stay text, so force
# them with a leading apostrophe the same way typing them in Excel would)
$ws.Cells.Item(5, 1).Value = "'1"
$ws.Cells.Item(5, 2).Value = "'23"
$ws.Cells.Item(5, 3).Value = "'43"
$ws.Cells.Item(5, 4).Value = "2025-01-30 11:49:51"
$ws.Cells.Item(5, 5).Value = "'65"
$ws.Cells.Item(5, 6).Value = "'76"
$ws.Cells.Item(5, 7).Value = "'54"
$ws.Cells.Item(5, 8).Value = "em dia"
$ws.Cells.Item(5, 9).Value = "em aberto"
$ws.Cells.Item(5, 10).Value = "'87"
$ws.Cells.Item(5, 11).Value = "sim"
